$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2032258064516129
$ws.Range("C2").Value = 0.5483870967741935
$ws.Range("J2").Value = 0.01612903225806452
$ws.Range("P2").Value = 0.1354838709677419
$ws.Range("S2").Value = 0.09677419354838709
# Row 3
$ws.Range("B3").Value = 0.005714285714285714
$ws.Range("C3").Value = 0.02285714285714286
$ws.Range("J3").Value = 0.06857142857142857
$ws.Range("P3").Value = 0.7428571428571429
$ws.Range("S3").Value = 0.16
# Row 4
$ws.Range("J4").Value = 0.1363636363636364
$ws.Range("P4").Value = 0.6818181818181818
$ws.Range("S4").Value = 0.1818181818181818
# Row 6
$ws.Range("B6").Value = 0.07692307692307693
$ws.Range("D6").Value = 0.02403846153846154
$ws.Range("F6").Value = 0.04807692307692308
$ws.Range("J6").Value = 0.2451923076923077
$ws.Range("O6").Value = 0.01923076923076923
$ws.Range("Q6").Value = 0.1682692307692308
$ws.Range("R6").Value = 0.1009615384615385
$ws.Range("S6").Value = 0.3173076923076923
# Row 7
$ws.Range("B7").Value = 0.08333333333333333
$ws.Range("D7").Value = 0.02976190476190476
$ws.Range("E7").Value = 0.005952380952380952
$ws.Range("F7").Value = 0.04761904761904762
$ws.Range("J7").Value = 0.1130952380952381
$ws.Range("O7").Value = 0.0119047619047619
$ws.Range("Q7").Value = 0.1785714285714286
$ws.Range("R7").Value = 0.07738095238095238
$ws.Range("S7").Value = 0.4523809523809524
# Row 8
$ws.Range("B8").Value = 0.1020833333333333
$ws.Range("D8").Value = 0.0375
$ws.Range("E8").Value = 0.002083333333333333
$ws.Range("F8").Value = 0.05625
$ws.Range("J8").Value = 0.1229166666666667
$ws.Range("O8").Value = 0.01041666666666667
$ws.Range("Q8").Value = 0.21875
$ws.Range("R8").Value = 0.09166666666666666
$ws.Range("S8").Value = 0.3583333333333333
# Row 9
$ws.Range("B9").Value = 0.1027027027027027
$ws.Range("D9").Value = 0.03783783783783784
$ws.Range("F9").Value = 0.02162162162162162
$ws.Range("J9").Value = 0.1675675675675676
$ws.Range("O9").Value = 0.01621621621621622
$ws.Range("Q9").Value = 0.2162162162162162
$ws.Range("R9").Value = 0.07567567567567568
$ws.Range("S9").Value = 0.3621621621621622
# Row 10
$ws.Range("B10").Value = 0.1066066066066066
$ws.Range("D10").Value = 0.02402402402402402
$ws.Range("E10").Value = 0.0007507507507507507
$ws.Range("F10").Value = 0.07507507507507508
$ws.Range("J10").Value = 0.1351351351351351
$ws.Range("O10").Value = 0.01276276276276276
$ws.Range("Q10").Value = 0.210960960960961
$ws.Range("R10").Value = 0.09984984984984985
$ws.Range("S10").Value = 0.3348348348348348
# Row 11
$ws.Range("G11").Value = 0.1124497991967871
$ws.Range("J11").Value = 0.1204819277108434
$ws.Range("K11").Value = 0.1686746987951807
$ws.Range("L11").Value = 0.5863453815261044
$ws.Range("S11").Value = 0.01204819277108434
# Row 12
$ws.Range("G12").Value = 0.7337662337662337
$ws.Range("J12").Value = 0.1948051948051948
$ws.Range("K12").Value = 0.01298701298701299
$ws.Range("L12").Value = 0.03246753246753246
$ws.Range("S12").Value = 0.02597402597402598
# Row 13
$ws.Range("G13").Value = 0.6818181818181818
$ws.Range("J13").Value = 0.2727272727272727
$ws.Range("S13").Value = 0.04545454545454546
# Row 15
$ws.Range("F15").Value = 0.02941176470588235
$ws.Range("H15").Value = 0.2058823529411765
$ws.Range("I15").Value = 0.05392156862745098
$ws.Range("J15").Value = 0.3137254901960784
$ws.Range("K15").Value = 0.05882352941176471
$ws.Range("M15").Value = 0.02450980392156863
$ws.Range("O15").Value = 0.06862745098039216
$ws.Range("S15").Value = 0.2450980392156863
# Row 16
$ws.Range("F16").Value = 0.02870813397129187
$ws.Range("H16").Value = 0.1531100478468899
$ws.Range("I16").Value = 0.09090909090909091
$ws.Range("J16").Value = 0.3827751196172249
$ws.Range("K16").Value = 0.1100478468899522
$ws.Range("M16").Value = 0.01913875598086124
$ws.Range("O16").Value = 0.04784688995215311
$ws.Range("S16").Value = 0.1674641148325359
# Row 17
$ws.Range("F17").Value = 0.01022494887525562
$ws.Range("H17").Value = 0.2208588957055215
$ws.Range("I17").Value = 0.08384458077709611
$ws.Range("J17").Value = 0.4028629856850716
$ws.Range("K17").Value = 0.09611451942740286
$ws.Range("M17").Value = 0.01431492842535787
$ws.Range("O17").Value = 0.05725971370143149
$ws.Range("S17").Value = 0.114519427402863
# Row 18
$ws.Range("F18").Value = 0.01351351351351351
$ws.Range("H18").Value = 0.2297297297297297
$ws.Range("I18").Value = 0.06756756756756757
$ws.Range("J18").Value = 0.4504504504504505
$ws.Range("K18").Value = 0.05405405405405406
$ws.Range("O18").Value = 0.06756756756756757
$ws.Range("S18").Value = 0.1171171171171171
# Row 19
$ws.Range("F19").Value = 0.0109519797809604
$ws.Range("H19").Value = 0.2097725358045493
$ws.Range("I19").Value = 0.08340353833192923
$ws.Range("J19").Value = 0.3934288121314238
$ws.Range("K19").Value = 0.08845829823083404
$ws.Range("M19").Value = 0.02358887952822241
$ws.Range("N19").Value = 0.003369839932603201
$ws.Range("O19").Value = 0.06908171861836562
$ws.Range("S19").Value = 0.117944397641112
